$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1424
$ws.Range("F5").Value = 6714
$ws.Range("F6").Value = 527
$ws.Range("F7").Value = 1058
$ws.Range("F8").Value = 38
$ws.Range("F9").Value = 4570
$ws.Range("F10").Value = 6786
$ws.Range("F12").Value = 225
$ws.Range("F13").Value = 1390
$ws.Range("F14").Value = 800
$ws.Range("F15").Value = 114
$ws.Range("F17").Value = 32
$ws.Range("F18").Value = 1129
$ws.Range("F20").Value = 128
$ws.Range("F22").Value = 185
$ws.Range("F24").Value = 1052
$ws.Range("F25").Value = 538
$ws.Range("F26").Value = 38
$ws.Range("F27").Value = 28
$ws.Range("F28").Value = 118
$ws.Range("F30").Value = 1164
$ws.Range("F31").Value = 30
$ws.Range("F32").Value = 96
$ws.Range("F34").Value = 4
$ws.Range("F38").Value = 363
$ws.Range("F39").Value = 39
$ws.Range("F41").Value = 312
$ws.Range("F42").Value = 1181
$ws.Range("F43").Value = 522
$ws.Range("F44").Value = 62
$ws.Range("F45").Value = 110

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 14
$ws.Range("F3").Value = 14
$ws.Range("F5").Value = 13
$ws.Range("F8").Value = 28
$ws.Range("F13").Value = 16
$ws.Range("F14").Value = 13
$ws.Range("F17").Value = 1728
$ws.Range("F30").Value = 110
$ws.Range("F31").Value = 766
$ws.Range("F32").Value = 962
$ws.Range("F33").Value = 581
$ws.Range("F35").Value = 90
$ws.Range("F38").Value = 96
$ws.Range("F39").Value = 128
$ws.Range("F41").Value = 2
$ws.Range("F42").Value = 60

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 714
$ws.Range("F5").Value = 842
$ws.Range("F6").Value = 610
$ws.Range("F8").Value = 1271
$ws.Range("F9").Value = 1769

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 14
$ws.Range("F3").Value = 714
$ws.Range("F4").Value = 1424
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 842
$ws.Range("F9").Value = 610
$ws.Range("F10").Value = 610
$ws.Range("F12").Value = 6714
$ws.Range("F13").Value = 527
$ws.Range("F14").Value = 1058
$ws.Range("F15").Value = 38
$ws.Range("F16").Value = 4570
$ws.Range("F18").Value = 6786
$ws.Range("F19").Value = 225
$ws.Range("F20").Value = 1390
$ws.Range("F22").Value = 800
$ws.Range("F23").Value = 114
$ws.Range("F24").Value = 1271
$ws.Range("F26").Value = 1129
$ws.Range("F27").Value = 128
$ws.Range("F28").Value = 185
$ws.Range("F29").Value = 1052
$ws.Range("F31").Value = 538
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 118
$ws.Range("F35").Value = 1164
$ws.Range("F36").Value = 96
$ws.Range("F39").Value = 962
$ws.Range("F41").Value = 581
$ws.Range("F42").Value = 363
$ws.Range("F43").Value = 39
$ws.Range("F44").Value = 90
$ws.Range("F45").Value = 312
$ws.Range("F46").Value = 522
$ws.Range("F47").Value = 96
$ws.Range("F49").Value = 110
$ws.Range("F50").Value = 60
